# Actualización automática 2025-12-01 08:30:07
#
# This workbook carries two sheets that get touched:
#   1) "VENTAS POR GRUPO"  -> a handful of now-stale product totals for this
#      advisor's 24-order window roll off to 0 (and the "x de 24" counters
#      for PIEDRA SINTERIZADA / PORCELANATO follow suit).
#   2) "VENTA MENSUAL"     -> the rolling 4-month window (agosto..noviembre)
#      advances one month (septiembre..diciembre): every client's monthly
#      figures shift left one column, the new right-most month (diciembre)
#      starts at 0 except for any new postings, and the column headers /
#      widths / column totals follow the same shift.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" - zero out orders that fell out of the window
# ---------------------------------------------------------------------
$ws1.Range("M4").Value  = 0
$ws1.Range("L6").Value  = 0
$ws1.Range("M6").Value  = 0
$ws1.Range("M12").Value = 0
$ws1.Range("M14").Value = 0
$ws1.Range("M25").Value = 0

# "x de 24" counters for PIEDRA SINTERIZADA (L) / PORCELANATO (M) drop too
$ws1.Range("L26").Value = "0 de 24"
$ws1.Range("M26").Value = "0 de 24"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" - roll the 4-month window forward by one month
# ---------------------------------------------------------------------

# Column widths follow the same column (C/D/E) shift; F/G stay put.
# NB: the COM ColumnWidth property and the raw OOXML <col width> attribute
# differ by a fixed +5/6 character padding offset in this engine, so the
# assigned values are pre-compensated to land exactly on 16 / 14 / 15.
$ws2.Columns.Item(3).ColumnWidth = 15.1666666666667
$ws2.Columns.Item(4).ColumnWidth = 13.1666666666667
$ws2.Columns.Item(5).ColumnWidth = 14.1666666666667

# Month headers shift left, December is newly introduced in column F
$ws2.Range("C1").Value = "septiembre"
$ws2.Range("D1").Value = "octubre"
$ws2.Range("E1").Value = "noviembre"
$ws2.Range("F1").Value = "diciembre"

# Per-client monthly figures: C<-D, D<-E, E<-F, F<-0, with the single
# exception of E22 which receives a genuinely new posting (5686.14)
# rather than the carried-over (zero) value.
$ws2.Range("C3").Value  = -3519.22
$ws2.Range("D3").Value  = 0
$ws2.Range("E3").Value  = 0
$ws2.Range("F3").Value  = 0

$ws2.Range("C4").Value  = 306.24
$ws2.Range("D4").Value  = 2548.88
$ws2.Range("E4").Value  = 1653.75
$ws2.Range("F4").Value  = 0

$ws2.Range("C5").Value  = 0
$ws2.Range("D5").Value  = 1362.43
$ws2.Range("E5").Value  = 0
$ws2.Range("F5").Value  = 0

$ws2.Range("C6").Value  = 1265.01
$ws2.Range("D6").Value  = 0
$ws2.Range("E6").Value  = 24929.52
$ws2.Range("F6").Value  = 0

$ws2.Range("C8").Value  = 2411.69
$ws2.Range("D8").Value  = 663.55
$ws2.Range("E8").Value  = 0
$ws2.Range("F8").Value  = 0

$ws2.Range("C9").Value  = 142.56
$ws2.Range("D9").Value  = 0
$ws2.Range("E9").Value  = 0
$ws2.Range("F9").Value  = 0

$ws2.Range("C10").Value = 549.5
$ws2.Range("D10").Value = 4220.84
$ws2.Range("E10").Value = 0
$ws2.Range("F10").Value = 0

$ws2.Range("C11").Value = -3989.12
$ws2.Range("D11").Value = 0
$ws2.Range("E11").Value = 0
$ws2.Range("F11").Value = 0

$ws2.Range("C12").Value = 17655.41
$ws2.Range("D12").Value = 6935.82
$ws2.Range("E12").Value = 5890.54
$ws2.Range("F12").Value = 0

$ws2.Range("C14").Value = 0
$ws2.Range("D14").Value = 0
$ws2.Range("E14").Value = 3080.12
$ws2.Range("F14").Value = 0

$ws2.Range("C16").Value = 829.4400000000001
$ws2.Range("D16").Value = 0
$ws2.Range("E16").Value = 0
$ws2.Range("F16").Value = 0

$ws2.Range("C17").Value = 0
$ws2.Range("D17").Value = 3162.93
$ws2.Range("E17").Value = 0
$ws2.Range("F17").Value = 0

$ws2.Range("C20").Value = 3252.41
$ws2.Range("D20").Value = 0
$ws2.Range("E20").Value = 0
$ws2.Range("F20").Value = 0

$ws2.Range("C21").Value = 1632.93
$ws2.Range("D21").Value = 1687.32
$ws2.Range("E21").Value = 0
$ws2.Range("F21").Value = 0

$ws2.Range("C22").Value = -347.92
$ws2.Range("D22").Value = 14679.01
$ws2.Range("E22").Value = 5686.14
$ws2.Range("F22").Value = 0

$ws2.Range("C25").Value = 3690.09
$ws2.Range("D25").Value = 6348.54
$ws2.Range("E25").Value = 23307.46
$ws2.Range("F25").Value = 0

# Column totals (row 26) follow the shifted data exactly
$ws2.Range("C26").Value = 23879.02
$ws2.Range("D26").Value = 41609.32
$ws2.Range("E26").Value = 64547.53
$ws2.Range("F26").Value = 0
